$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.82%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'27.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.54%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'4.864"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.26%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.06414"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.39%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'6.940"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.96%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.183"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-5.84%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.8757"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.73%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.51%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.05133"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.22%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07501"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.54%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.02943"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.35%"
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'-0.42%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001567"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.59%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0006402"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.15%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.006118"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'4.56%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D18").Value = "'3.305"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.56%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'0.16%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D22").Value = "'3.904"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.03%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04424"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.72%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D25").Value = "'0.001177"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.09%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'-8.95%"
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'8.24%"
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'15.04%"
$ws.Range("E28").Style = "Normal"

$ws.Range("D40").Value = "'0.04174"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.21%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.006797"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.34%"
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'0.66%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D44").Value = "'0.01147"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'7.12%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005302"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.10%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.685"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'13.33%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'-7.40%"
$ws.Range("E47").Style = "Normal"
